# Apply the "working script to plot sample locations" edit:
#  - insert a new "meas_type" column (C) between "stick" and "x_lo"
#  - populate meas_type values for every existing row
#  - fix up a few x_lo/x_hi/y_lo/y_hi values that changed alongside the
#    column insert (rows 7 and 10)
#  - append two new sample rows (ar / thin_section_vert)
#  - move the active selection to E11

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new blank column at C (shifts old C:H -> D:I).
$ws.Columns("C:C").Insert()

# 2) Row 1 header.
$ws.Range("A1").Value = 'section'
$ws.Range("B1").Value = 'stick'
$ws.Range("C1").Value = 'meas_type'
$ws.Range("D1").Value = 'x_lo'
$ws.Range("E1").Value = 'x_hi'
$ws.Range("F1").Value = 'y_lo'
$ws.Range("G1").Value = 'y_hi'
$ws.Range("H1").Value = 'section_top_depth'
$ws.Range("I1").Value = 'offset'

# 3) Data rows 2-15: fill in the new meas_type column, and correct the
#    handful of x/y values that changed.
$ws.Range("C2").Value = 'IC / CC'
$ws.Range("C3").Value = 'IC / CC'
$ws.Range("C4").Value = 'ICMPS'
$ws.Range("C5").Value = 'CO$_2$'
$ws.Range("C6").Value = 'CH$_4$'

$ws.Range("C7").Value = 'ISO'
$ws.Range("F7").Value = -110
$ws.Range("G7").Value = -120

$ws.Range("C8").Value = 'ISO'
$ws.Range("C9").Value = 'ISO'

$ws.Range("C10").Value = 'ISO'
$ws.Range("E10").Value = 120

$ws.Range("C11").Value = 'ISO'
$ws.Range("C12").Value = 'CC'
$ws.Range("C13").Value = 'ICPMS'
$ws.Range("C14").Value = 'IC'
$ws.Range("C15").Value = 'CO$_2$'

# Uniform formatting for the whole new meas_type column (matches the
# look of the other populated cells in the sheet).
$ws.Range("C1:C17").Font.Name = "Aptos Narrow"
$ws.Range("C1:C17").Font.Size = 12

# 4) Append the two new sample rows.
$ws.Range("A16").Value = '228_4'
$ws.Range("B16").Value = 'ar'
$ws.Range("C16").Value = 'Age'
$ws.Range("D16").Value = -120
$ws.Range("E16").Value = 120
$ws.Range("F16").Value = -110
$ws.Range("G16").Value = -10

$ws.Range("A17").Value = '228_4'
$ws.Range("B17").Value = 'thin_section_vert'
$ws.Range("C17").Value = 'PP'
$ws.Range("D17").Value = -120
$ws.Range("E17").Value = 120
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = -10

$ws.Range("A16:G17").Font.Name = "Aptos Narrow"
$ws.Range("A16:G17").Font.Size = 12

# 5) Selection, matching the saved cursor position.
$ws.Range("E11").Select()
